$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (plain text columns) ---
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

# --- Price column (D): force text storage so values like "443.60" / "0.0690" keep exact formatting ---
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.688.81'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.355.16'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '193.57'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '592.85'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.76'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.424'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.935.28'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.48'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '69.676.80'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.368.43'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.84'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '443.60'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.75'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.79'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.65'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.489.75'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.518'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.193'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.58'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '23.13'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.62'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.28'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.05'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '164.63'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '27.25'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.813'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.60'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.764.37'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.53'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '25.52'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0690'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '344.84'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '40.61'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '32.59'

# --- Volume(1h) column (E): percentage text values ---
$ws.Range('E2').Value = '  +3.35%  '
$ws.Range('E3').Value = '  +4.34%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +5.84%  '
$ws.Range('E6').Value = '  +2.55%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('E9').Value = '  +3.84%  '
$ws.Range('E10').Value = '  +2.73%  '
$ws.Range('E11').Value = '  +2.80%  '
$ws.Range('E12').Value = '  +4.19%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('E15').Value = '  +3.25%  '
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('E19').Value = '  +13.09%  '
$ws.Range('E20').Value = '  +2.72%  '
$ws.Range('E21').Value = '  +3.76%  '
$ws.Range('E22').Value = '  +4.30%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  +3.98%  '
$ws.Range('E25').Value = '  +4.34%  '
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('E27').Value = '  +3.87%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('E30').Value = '  +2.95%  '
$ws.Range('E31').Value = '  +2.62%  '
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  +2.01%  '
$ws.Range('E37').Value = '  +3.07%  '
$ws.Range('E38').Value = '  +2.55%  '
$ws.Range('E39').Value = '  +4.18%  '
$ws.Range('E40').Value = '  +1.48%  '
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('E42').Value = '  +6.37%  '
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('E44').Value = '  +3.51%  '
$ws.Range('E45').Value = '  +3.68%  '
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('E47').Value = '  +3.43%  '
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('E50').Value = '  +5.60%  '
$ws.Range('E51').Value = '  +5.14%  '
